# Auto-generated edit script applying scheduled market-data refresh
# to the Sargatanas_Profits workbook (currentAveragePrice / Leve profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86
$ws.Range("H86").Value = 231482830
$ws.Range("J86").Value = 27778928
$ws.Range("L86").Value = 27778928
$ws.Range("N86").Value = -27781174
# Row 89
$ws.Range("H89").Value = 231482830
$ws.Range("J89").Value = 27778928
$ws.Range("L89").Value = 138894640
$ws.Range("N89").Value = -138905872
# Row 103
$ws.Range("H103").Value = 379.33334
$ws.Range("I103").Value = 369
$ws.Range("J103").Value = 400
$ws.Range("K103").Value = 1107
$ws.Range("L103").Value = 1200
$ws.Range("M103").Value = -521
$ws.Range("N103").Value = -2372
# Row 106
$ws.Range("H106").Value = 1754.2609
$ws.Range("I106").Value = 1754.2609
$ws.Range("K106").Value = 1754.2609
$ws.Range("M106").Value = -1123.2609
# Row 132
$ws.Range("H132").Value = 2229.7144
$ws.Range("I132").Value = 1851.4667
$ws.Range("K132").Value = 5554.4001
$ws.Range("M132").Value = -3024.4001
# Row 137
$ws.Range("H137").Value = 4406.1665
$ws.Range("I137").Value = 5880
$ws.Range("J137").Value = 2932.3333
$ws.Range("K137").Value = 17640
$ws.Range("L137").Value = 8796.999899999999
$ws.Range("M137").Value = -15090
$ws.Range("N137").Value = -13896.9999
# Row 138
$ws.Range("H138").Value = 1476051.8
$ws.Range("J138").Value = 2047120.6
$ws.Range("L138").Value = 6141361.800000001
$ws.Range("N138").Value = -6151641.800000001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3141.4697
$ws.Range("I32").Value = 3068.8413
$ws.Range("K32").Value = 3068.8413
$ws.Range("M32").Value = -2781.8413
# Row 61
$ws.Range("H61").Value = 5181.1904
$ws.Range("I61").Value = 2928.2693
$ws.Range("J61").Value = 15831.363
$ws.Range("K61").Value = 2928.2693
$ws.Range("L61").Value = 15831.363
$ws.Range("M61").Value = -2716.2693
$ws.Range("N61").Value = -16255.363
# Row 109
$ws.Range("H109").Value = 57586.332
$ws.Range("J109").Value = 57586.332
$ws.Range("L109").Value = 57586.332
$ws.Range("N109").Value = -60360.332
# Row 110
$ws.Range("H110").Value = 20836426
$ws.Range("I110").Value = 2597.4
$ws.Range("J110").Value = 55559476
$ws.Range("K110").Value = 2597.4
$ws.Range("L110").Value = 55559476
$ws.Range("M110").Value = -552.4000000000001
$ws.Range("N110").Value = -55563566
# Row 136
$ws.Range("H136").Value = 5181.1904
$ws.Range("I136").Value = 2928.2693
$ws.Range("J136").Value = 15831.363
$ws.Range("K136").Value = 8784.8079
$ws.Range("L136").Value = 47494.089
$ws.Range("M136").Value = -6234.8079
$ws.Range("N136").Value = -52594.089

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 20837202
$ws.Range("I20").Value = 27780918
$ws.Range("J20").Value = 6049
$ws.Range("K20").Value = 27780918
$ws.Range("L20").Value = 6049
$ws.Range("M20").Value = -27780671
$ws.Range("N20").Value = -6543
# Row 22
$ws.Range("H22").Value = 9490.091
$ws.Range("I22").Value = 9490.091
$ws.Range("K22").Value = 9490.091
$ws.Range("M22").Value = -9317.091

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 171.26086
$ws.Range("I7").Value = 108.52941
$ws.Range("K7").Value = 108.52941
$ws.Range("M7").Value = 4.470590000000001
# Row 22
$ws.Range("H22").Value = 373.05264
$ws.Range("I22").Value = 299.29413
$ws.Range("K22").Value = 299.29413
$ws.Range("M22").Value = 50.70587
# Row 31
$ws.Range("H31").Value = 5599.369
$ws.Range("J31").Value = 7448.189
$ws.Range("L31").Value = 7448.189
$ws.Range("N31").Value = -8038.189
# Row 34
$ws.Range("H34").Value = 5599.369
$ws.Range("J34").Value = 7448.189
$ws.Range("L34").Value = 7448.189
$ws.Range("N34").Value = -7852.189
# Row 62
$ws.Range("H62").Value = 7482.3335
$ws.Range("I62").Value = 6968
$ws.Range("J62").Value = 7739.5
$ws.Range("K62").Value = 6968
$ws.Range("L62").Value = 7739.5
$ws.Range("M62").Value = -6344
$ws.Range("N62").Value = -8987.5
# Row 65
$ws.Range("H65").Value = 7482.3335
$ws.Range("I65").Value = 6968
$ws.Range("J65").Value = 7739.5
$ws.Range("K65").Value = 34840
$ws.Range("L65").Value = 38697.5
$ws.Range("M65").Value = -31720
$ws.Range("N65").Value = -44937.5
# Row 99
$ws.Range("H99").Value = 9039.388999999999
$ws.Range("I99").Value = 9385
$ws.Range("K99").Value = 9385
$ws.Range("M99").Value = -7887
# Row 126
$ws.Range("H126").Value = 9039.388999999999
$ws.Range("I126").Value = 9385
$ws.Range("K126").Value = 28155
$ws.Range("M126").Value = -25685
# Row 132
$ws.Range("H132").Value = 3729.25
$ws.Range("I132").Value = 2034.4286
$ws.Range("K132").Value = 6103.2858
$ws.Range("M132").Value = -3573.2858
# Row 134
$ws.Range("H134").Value = 6333.3887
$ws.Range("I134").Value = 2625.0557
$ws.Range("J134").Value = 10041.723
$ws.Range("K134").Value = 7875.1671
$ws.Range("L134").Value = 30125.169
$ws.Range("M134").Value = -5340.1671
$ws.Range("N134").Value = -35195.169

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 3182
$ws.Range("J68").Value = 3368.2693
$ws.Range("L68").Value = 10104.8079
$ws.Range("N68").Value = -11726.8079
# Row 71
$ws.Range("H71").Value = 3182
$ws.Range("J71").Value = 3368.2693
$ws.Range("L71").Value = 30314.4237
$ws.Range("N71").Value = -38426.4237
# Row 132
$ws.Range("H132").Value = 7426.0303
$ws.Range("J132").Value = 7950.4165
$ws.Range("L132").Value = 71553.7485
$ws.Range("N132").Value = -76613.7485
# Row 140
$ws.Range("H140").Value = 288167
$ws.Range("I140").Value = 446104.22
$ws.Range("K140").Value = 1338312.66
$ws.Range("M140").Value = -1333132.66

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 4000101.2
$ws.Range("I2").Value = 135.66667
$ws.Range("J2").Value = 10000050
$ws.Range("K2").Value = 135.66667
$ws.Range("L2").Value = 10000050
$ws.Range("M2").Value = -22.66667000000001
$ws.Range("N2").Value = -10000276
# Row 132
$ws.Range("H132").Value = 8156.6
$ws.Range("I132").Value = 5590.5
$ws.Range("J132").Value = 13288.8
$ws.Range("K132").Value = 16771.5
$ws.Range("L132").Value = 39866.39999999999
$ws.Range("M132").Value = -14241.5
$ws.Range("N132").Value = -44926.39999999999
# Row 134
$ws.Range("H134").Value = 80256
$ws.Range("J134").Value = 80256
$ws.Range("L134").Value = 240768
$ws.Range("N134").Value = -245838

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1920564.4
$ws.Range("J46").Value = 5063.5
$ws.Range("L46").Value = 5063.5
$ws.Range("N46").Value = -5439.5
# Row 132
$ws.Range("H132").Value = 12201776
$ws.Range("I132").Value = 25003466
$ws.Range("K132").Value = 75010398
$ws.Range("M132").Value = -75007868

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 5394.1113
$ws.Range("J62").Value = 5825.25
$ws.Range("L62").Value = 5825.25
$ws.Range("N62").Value = -7073.25
# Row 65
$ws.Range("H65").Value = 5394.1113
$ws.Range("J65").Value = 5825.25
$ws.Range("L65").Value = 29126.25
$ws.Range("N65").Value = -35366.25
# Row 69
$ws.Range("H69").Value = 271
$ws.Range("J69").Value = 271
$ws.Range("L69").Value = 271
$ws.Range("N69").Value = -1769
# Row 72
$ws.Range("H72").Value = 271
$ws.Range("J72").Value = 271
$ws.Range("L72").Value = 813
$ws.Range("N72").Value = -8301
# Row 132
$ws.Range("H132").Value = 8630826
$ws.Range("I132").Value = 11370234
$ws.Range("J132").Value = 21258.428
$ws.Range("K132").Value = 34110702
$ws.Range("L132").Value = 63775.284
$ws.Range("M132").Value = -34108172
$ws.Range("N132").Value = -68835.284
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 136
$ws.Range("H136").Value = 37078244
$ws.Range("I136").Value = 55557812
$ws.Range("K136").Value = 166673436
$ws.Range("M136").Value = -166670886
# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
